# Updates the cryptos list in Sheet1 to reflect the latest scraped
# coinranking.com snapshot: refreshed Price (column D) and Volume(1h)
# (column E) figures for most rows, plus a rank swap between the
# Chainlink and WrappedEther rows (19/20), whose Coin, Link, Price and
# Volume values now trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.595.35'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.956.35'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.54'
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.46'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.953.67'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  -3.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  -4.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -2.36%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.609.45'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.442.69'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.04'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.949.41'
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.94'
$ws.Range("E20").Value = '  +14.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.93'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.695'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.06'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  -3.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.29'
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("E27").Value = '  -5.04%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  +8.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.12'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("E33").Value = '  +4.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.18'
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.72'
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.11'
$ws.Range("E38").Value = '  +5.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.08'
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("E40").Value = '  -7.06%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.303'
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.82'
$ws.Range("E43").Value = '  -6.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.53'
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '384.85'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0351'
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.681.33'
$ws.Range("E47").Value = '  -4.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.91'
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.90'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("E51").Value = '  +1.15%  '